$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.949.56'
$ws.Range('E2').Value = '  +2.93%  '
$ws.Range('D3').Value = '2.446.33'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.24%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.173'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.37%  '
$ws.Range('D10').Value = '2.447.05'
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.335'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.92%  '
$ws.Range('E13').Value = '  -1.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000180'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.81%  '
$ws.Range('D15').Value = '69.855.22'
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('D16').Value = '2.901.39'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.88%  '
$ws.Range('D18').Value = '2.444.24'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.85'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '341.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.88%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.38'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.25%  '
$ws.Range('D27').Value = '2.577.75'
$ws.Range('E27').Value = '  +0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.36%  '
$ws.Range('D30').Value = '0.0₃0858'
$ws.Range('E30').Value = '  +6.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.39'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '461.41'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.20%  '
$ws.Range('E33').Value = '  +10.18%  '
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('E37').Value = '  +6.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.09'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.23'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.303'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.10%  '
$ws.Range('E42').Value = '  +4.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.42'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '134.47'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('E49').Value = '  +2.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.491'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.564'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.55%  '
